$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update K2 text value
$ws.Range("K2").Value = "van der Weijst and Peterse (Unpublished data)"

# Update M2 and N2 numeric values
$ws.Range("M2").Value = 27.31
$ws.Range("N2").Value = 0.2875847710503443

# Remove columns O and P entirely (headers O1/P1 and data O2/P2)
$ws.Range("O1:P1").EntireColumn.Delete()
